$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 3
Write-Output ("ScrollColumn=" + $excel.ActiveWindow.ScrollColumn)
Write-Output ("ScrollRow=" + $excel.ActiveWindow.ScrollRow)
